$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Give the first paragraph a thin box border (5-twip gap on every side,
# no explicit line weight) and bump its left indent from 120 to 225
# twips (6pt -> 11.25pt) to match the bordered paragraphs later in the
# document.
$p1.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromRight = 5
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Locate the placeholder merge-field id text and swap it for the
# AFICC-specific id.
$findRange = $d.Content.Duplicate
$null = $findRange.Find.Execute("**ID__AFFARS_pgi_5301_topic_20__ID**", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$idStart = $findRange.Start
$idEnd = $findRange.End

$newText = "**ID__AFFARS_AFICC_PGI_5301_9001__ID**"
$idRange = $d.Range($idStart, $idEnd)
$idRange.Text = $newText

# The old text was immediately followed by a run containing a single
# trailing space; remove that run entirely now that the id has been
# replaced in place.
$spaceStart = $idStart + $newText.Length
$spaceRange = $d.Range($spaceStart, $spaceStart + 1)
$spaceRange.Delete()
